$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# 1. Remove the old "Control Data Size" mini-table (rows 13-14 of the
#    original layout): a 2-row header/data block that is being replaced by
#    the richer BFAST/BWA/BOWTIE comparison table below.
# ---------------------------------------------------------------------------
$ws.Rows(13).Delete()
$ws.Rows(13).Delete()

# ---------------------------------------------------------------------------
# 2. Make room for the new BFAST(400MB)/BWA/BOWTIE rows by pushing the
#    "App Temp for mapping" block (and everything under it) down by 4 rows.
# ---------------------------------------------------------------------------
$ws.Range("A14:L17").Insert()

# ---------------------------------------------------------------------------
# 3. Tool names for the two new comparison blocks.
# ---------------------------------------------------------------------------
$ws.Range("F10").Value = "BWA"
$ws.Range("F11").Value = "BWA"
$ws.Range("F12").Value = "BWA"
$ws.Range("F14").Value = "BOWTIE"
$ws.Range("F15").Value = "BOWTIE"
$ws.Range("F16").Value = "BOWTIE"

# ---------------------------------------------------------------------------
# 4. Third BFAST row (400MB treat/control size) directly under the existing
#    BFAST 100MB/200MB rows.
# ---------------------------------------------------------------------------
$ws.Range("B8").Value = "400MB"
$ws.Range("C8").Value = "400MB"
$ws.Range("D8").Value = "mouse ch19"
$ws.Range("E8").Value = "3.1G"
$ws.Range("F8").Value = "BFAST"

# ---------------------------------------------------------------------------
# 5. BWA block (rows 10-12): Treat/Control size, genome, mapping time, and
#    the resulting SAM file sizes.
# ---------------------------------------------------------------------------
$ws.Range("B10").Value = "100MB"
$ws.Range("C10").Value = "100MB"
$ws.Range("D10").Value = "mouse ch19"
$ws.Range("E10").Value = "3.1G"
$ws.Range("G10").Value = 182.96
$ws.Range("H10").Value = "114 MB"
$ws.Range("I10").Value = "122 MB"

$ws.Range("B11").Value = "200MB"
$ws.Range("C11").Value = "200MB"
$ws.Range("D11").Value = "mouse ch19"
$ws.Range("E11").Value = "3.1G"
$ws.Range("G11").Value = 307.47

$ws.Range("B12").Value = "400MB"
$ws.Range("C12").Value = "400MB"
$ws.Range("D12").Value = "mouse ch19"
$ws.Range("E12").Value = "3.1G"
$ws.Range("G12").Value = 542.65

$ws.Range("I11").Value = "245MB"
$ws.Range("I12").Value = "480 MB"
$ws.Range("H12").Value = "449 MB"
$ws.Range("H11").Value = "229 MB"

# ---------------------------------------------------------------------------
# 6. BOWTIE block (rows 14-16): same layout as the BWA block.
# ---------------------------------------------------------------------------
$ws.Range("B14").Value = "100MB"
$ws.Range("C14").Value = "100MB"
$ws.Range("D14").Value = "mouse ch19"
$ws.Range("E14").Value = "3.1G"
$ws.Range("G14").Value = 167.18

$ws.Range("B15").Value = "200MB"
$ws.Range("C15").Value = "200MB"
$ws.Range("D15").Value = "mouse ch19"
$ws.Range("E15").Value = "3.1G"
$ws.Range("G15").Value = 277.59

$ws.Range("B16").Value = "400MB"
$ws.Range("C16").Value = "400MB"
$ws.Range("D16").Value = "mouse ch19"
$ws.Range("E16").Value = "3.1G"
$ws.Range("G16").Value = 507.61

$ws.Range("H16").Value = "23MB"
$ws.Range("I16").Value = "23MB"
$ws.Range("H15").Value = "12Mb"
$ws.Range("I15").Value = "12Mb"
$ws.Range("H14").Value = "5.7 MB"
$ws.Range("I14").Value = "5.7 MB"

# ---------------------------------------------------------------------------
# 7. Row 16 (the 400MB BOWTIE row) got highlighted with a black font color
#    by the author -- apply the same formatting, matching the original
#    cell-by-cell selection (the Tool cell, F16, was left untouched).
# ---------------------------------------------------------------------------
$ws.Range("A16:E16").Font.Color = 0
$ws.Range("G16:L16").Font.Color = 0

# ---------------------------------------------------------------------------
# 8. Restore the active-cell selection left behind by the author.
# ---------------------------------------------------------------------------
$ws.Range("I15").Select()
